# Split the code line "type T is array[10] of Integer;" into three runs so
# that it reads "type T = array[10] of Integer;" (the "is" becomes "=",
# each of the three resulting pieces being its own run, matching how the
# slide was hand-edited in PowerPoint).
#
#   Run 1: "type "
#   Run 2: "T = "                  (replaces the original "T is ")
#   Run 3: "array[10] of Integer;"

$p = $ppt.ActivePresentation

$targetText = "type T is array[10] of Integer;"
$found = $false

for ($si = 1; $si -le $p.Slides.Count -and -not $found; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count -and -not $found; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }

        $tr = $shape.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count

        for ($pi = 1; $pi -le $paraCount -and -not $found; $pi++) {
            $para = $tr.Paragraphs($pi, 1)

            if ($para.Text.StartsWith($targetText)) {
                # Characters() is 1-based over the paragraph's own text.
                #   "type T is array[10] of Integer;"
                #    1234567890123456789012345678901
                # "type "  -> chars 1-5   (own run, unchanged text)
                # "T is "  -> chars 6-10  (own run, text becomes "T = ")
                # "array[10] of Integer;" -> chars 11-31 (own run, unchanged)

                $run1 = $para.Characters(1, 5)
                $run1.Text = "type "

                $run2 = $para.Characters(6, 5)
                $run2.Text = "T = "

                $run3 = $para.Characters(10, 22)
                $run3.Text = "array[10] of Integer;"

                $found = $true
            }
        }
    }
}
